$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -2
$ws.Range("F3").Value = -1
$ws.Range("F7").Value = -6
$ws.Range("F8").Value = -7
$ws.Range("F9").Value = -6
$ws.Range("F10").Value = -7
$ws.Range("F11").Value = -1
$ws.Range("F36").Value = 7
$ws.Range("F41").Value = 0
$ws.Range("F42").Value = 0
